$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (row, date-serial, col B, col C, col D)
$data = @(
    @(375, 44449, 1, 4, 66.54466810846782),
    @(376, 44450, 0, 4, 66.54466810846782),
    @(377, 44451, 1, 4, 66.54466810846782),
    @(378, 44452, 1, 5, 83.18083513558476),
    @(379, 44453, 1, 5, 83.18083513558476),
    @(380, 44454, 0, 5, 83.18083513558476),
    @(381, 44455, 0, 4, 66.54466810846782),
    @(382, 44456, 2, 5, 83.18083513558476),
    @(383, 44457, 0, 5, 83.18083513558476),
    @(384, 44458, 0, 4, 66.54466810846782),
    @(385, 44459, 0, 3, 49.90850108135086)
)

foreach ($row in $data) {
    $r = $row[0]

    # Copy formatting (style) from the row above into column A of the new row
    $ws.Range("A" + ($r - 1)).Copy($ws.Range("A" + $r))

    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
